$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: E20 10 -> 9 ; F20 gets a new grading comment
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "(-1)For incorrectly getting products of a customer."

# Row 34: E34 0 -> 3 ; F34 comment replaced, wrap text applied, row height grows to fit the longer comment
$ws.Range("E34").Value = 3
$ws.Range("F34").Value = "(-4)I have changed your addProduct() code and run the test cases then 4 test cases failed but I didn’t deducted any points for remaining test cases`n"
$ws.Range("F34").WrapText = $true
$ws.Rows(34).RowHeight = 75

# Move the active selection to G34, matching where the grader was last working
$ws.Range("G34").Select() | Out-Null
